$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.514.72"
$ws.Range("E2").Value = "  -3.50%  "

$ws.Range("D3").Value = "1.848.70"
$ws.Range("E3").Value = "  -3.75%  "

$ws.Range("E4").Value = "  -0.89%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.74"
$ws.Range("E5").Value = "  +2.27%  "

$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4666"
$ws.Range("E7").Value = "  -3.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3922"
$ws.Range("E8").Value = "  -3.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.46"
$ws.Range("E9").Value = "  -2.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07909"
$ws.Range("E10").Value = "  -4.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9846"

$ws.Range("E12").Value = "  -5.36%  "

$ws.Range("D13").Value = "1.974.13"
$ws.Range("E13").Value = "  +2.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.846"
$ws.Range("E14").Value = "  -3.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.006"
$ws.Range("E15").Value = "  -3.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06855"
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.71"
$ws.Range("E18").Value = "  -4.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001006"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.06"
$ws.Range("E20").Value = "  -3.07%  "

$ws.Range("E21").Value = "  -0.79%  "

$ws.Range("D22").Value = "28.543.27"
$ws.Range("E22").Value = "  -3.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.385"
$ws.Range("E23").Value = "  -5.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  -5.46%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.127"
$ws.Range("E25").Value = "  -2.70%  "

$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.165.60"
$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.60"
$ws.Range("E27").Value = "  -1.43%  "

$ws.Range("E28").Value = "  -3.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.130"
$ws.Range("E29").Value = "  -6.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.015"
$ws.Range("E30").Value = "  -4.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.41"
$ws.Range("E31").Value = "  -2.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9776"
$ws.Range("E32").Value = "  -4.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09418"
$ws.Range("E33").Value = "  -2.28%  "

$ws.Range("E34").Value = "  -4.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.500"
$ws.Range("E35").Value = "  -1.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.345"
$ws.Range("E36").Value = "  -2.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06124"
$ws.Range("E37").Value = "  -3.69%  "

$ws.Range("E38").Value = "  -4.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.160"
$ws.Range("E39").Value = "  -2.64%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5697"
$ws.Range("E40").Value = "  -4.35%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.591"
$ws.Range("E41").Value = "  -4.21%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.09"
$ws.Range("E42").Value = "  -6.26%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1792"
$ws.Range("E43").Value = "  -3.17%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.394"
$ws.Range("E44").Value = "  -3.39%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.224"
$ws.Range("E45").Value = "  -4.76%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.80"
$ws.Range("E46").Value = "  -4.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5370"
$ws.Range("E47").Value = "  -3.63%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07157"
$ws.Range("E48").Value = "  -4.52%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.901"
$ws.Range("E49").Value = "  -2.42%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.57"
$ws.Range("E50").Value = "  -4.41%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.85"
$ws.Range("E51").Value = "  +1.53%  "
